# Loan RBI, Variable Instalments
# Insert a new (blank) column before column N on the "Repayment schedule"
# sheet, pushing the existing "Late" / heading / "Outstanding" columns one
# place to the right, then switch to that sheet with cell R6 selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Inherit the new column's width from the column immediately to its left,
# the way Excel's own "Insert Column" behaves.
$leftWidth = $ws.Columns("M").ColumnWidth

$ws.Columns("N").Insert()
$ws.Columns("N").ColumnWidth = $leftWidth

# Make "Repayment schedule" the active sheet/tab with R6 selected (this
# also clears tabSelected on whichever sheet was active before).
$ws.Activate()
$ws.Range("R6").Select() | Out-Null
